# Auto update Excel log
# Appends newly logged sensor events to the PIR, Humidity, Proximity and
# Camera sheets of the SeniorConnect master log workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# PIR sheet: new rows 14-26 (Bathroom / No Motion / Inactive)
# ---------------------------------------------------------------------
$pirRows = @(
    @("2026-01-30","16:00:23","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","16:00:24","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","16:00:29","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","16:00:34","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","16:00:39","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","16:00:44","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","16:00:49","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","16:00:54","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","16:00:59","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","16:01:04","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","16:01:09","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","16:01:14","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","16:01:19","16:00","Bathroom","No Motion","Inactive")
)

$ws = $wb.Worksheets.Item("PIR")
$r = 14
foreach ($row in $pirRows) {
    $ws.Cells.Item($r,1).NumberFormat = "@"
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,5).Value = $row[4]
    $ws.Cells.Item($r,6).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Humidity sheet: new rows 10-20 (Bathroom / percentage / Active)
# ---------------------------------------------------------------------
$humidityRows = @(
    @("2026-01-30","16:00:23","16:00","Bathroom","86.7%","Active"),
    @("2026-01-30","16:00:24","16:00","Bathroom","87.7%","Active"),
    @("2026-01-30","16:00:29","16:00","Bathroom","87.7%","Active"),
    @("2026-01-30","16:00:34","16:00","Bathroom","87.6%","Active"),
    @("2026-01-30","16:00:44","16:00","Bathroom","87.7%","Active"),
    @("2026-01-30","16:00:49","16:00","Bathroom","87.6%","Active"),
    @("2026-01-30","16:00:54","16:00","Bathroom","87.6%","Active"),
    @("2026-01-30","16:01:00","16:00","Bathroom","86.7%","Active"),
    @("2026-01-30","16:01:04","16:00","Bathroom","87.6%","Active"),
    @("2026-01-30","16:01:09","16:00","Bathroom","86.7%","Active"),
    @("2026-01-30","16:01:19","16:00","Bathroom","86.7%","Active")
)

$ws = $wb.Worksheets.Item("Humidity")
$r = 10
foreach ($row in $humidityRows) {
    $ws.Cells.Item($r,1).NumberFormat = "@"
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,5).NumberFormat = "@"
    $ws.Cells.Item($r,5).Value = $row[4]
    $ws.Cells.Item($r,6).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Proximity sheet: new rows 5-8 (Living Room Main Door ENTER/EXIT)
# ---------------------------------------------------------------------
$proximityRows = @(
    @("2026-01-30","16:00:37","16:00","Living Room Main Door","EXIT","User EXITED Living Room Main Door"),
    @("2026-01-30","16:00:39","16:00","Living Room Main Door","ENTER","User ENTERED Living Room Main Door"),
    @("2026-01-30","16:01:01","16:00","Living Room Main Door","EXIT","User EXITED Living Room Main Door"),
    @("2026-01-30","16:01:03","16:00","Living Room Main Door","ENTER","User ENTERED Living Room Main Door")
)

$ws = $wb.Worksheets.Item("Proximity")
$r = 5
foreach ($row in $proximityRows) {
    $ws.Cells.Item($r,1).NumberFormat = "@"
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,5).Value = $row[4]
    $ws.Cells.Item($r,6).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Camera sheet: new rows 5-8 (Living Room Main Door image captures)
# ---------------------------------------------------------------------
$cameraRows = @(
    @("2026-01-30","16:00:37","16:00","Living Room Main Door","Image Captured (EXIT)","Active"),
    @("2026-01-30","16:00:39","16:00","Living Room Main Door","Image Captured (ENTER)","Active"),
    @("2026-01-30","16:01:01","16:00","Living Room Main Door","Image Captured (EXIT)","Active"),
    @("2026-01-30","16:01:03","16:00","Living Room Main Door","Image Captured (ENTER)","Active")
)

$ws = $wb.Worksheets.Item("Camera")
$r = 5
foreach ($row in $cameraRows) {
    $ws.Cells.Item($r,1).NumberFormat = "@"
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,5).Value = $row[4]
    $ws.Cells.Item($r,6).Value = $row[5]
    $r = $r + 1
}
